$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "65.951.19"
$ws.Range("E2").Value = "  +1.72%  "

$ws.Range("D3").Value = "3.212.42"
$ws.Range("E3").Value = "  +1.39%  "

$ws.Range("E4").Value = "  -0.07%  "

$ws.Range("D5").Value = "'606.09"
$ws.Range("E5").Value = "  +4.66%  "

$ws.Range("D6").Value = "'153.35"
$ws.Range("E6").Value = "  +1.28%  "

$ws.Range("E7").Value = "  -0.14%  "

$ws.Range("D8").Value = "3.211.74"
$ws.Range("E8").Value = "  +1.49%  "

$ws.Range("D9").Value = "'0.533"
$ws.Range("E9").Value = "  +0.12%  "

$ws.Range("E10").Value = "  -0.73%  "

$ws.Range("D11").Value = "'6.16"
$ws.Range("E11").Value = "  -1.11%  "

$ws.Range("D12").Value = "'0.509"
$ws.Range("E12").Value = "  +1.51%  "

$ws.Range("E13").Value = "  +1.09%  "

$ws.Range("D14").Value = "'39.08"
$ws.Range("E14").Value = "  +3.82%  "

$ws.Range("D15").Value = "3.737.45"
$ws.Range("E15").Value = "  +1.31%  "

$ws.Range("B16").Value = "Polkadot"
$ws.Range("C16").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D16").Value = "'7.46"
$ws.Range("E16").Value = "  +4.07%  "

$ws.Range("B17").Value = "WrappedBTC"
$ws.Range("C17").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D17").Value = "66.104.90"
$ws.Range("E17").Value = "  +1.71%  "

$ws.Range("D18").Value = "3.204.86"
$ws.Range("E18").Value = "  +1.25%  "

$ws.Range("E19").Value = "  -0.18%  "

$ws.Range("D20").Value = "'509.89"
$ws.Range("E20").Value = "  -0.40%  "

$ws.Range("E21").Value = "  +4.62%  "

$ws.Range("D22").Value = "'0.734"
$ws.Range("E22").Value = "  +0.80%  "

$ws.Range("D23").Value = "'15.31"
$ws.Range("E23").Value = "  +0.36%  "

$ws.Range("D24").Value = "'8.05"
$ws.Range("E24").Value = "  +3.40%  "

$ws.Range("D25").Value = "'85.26"
$ws.Range("E25").Value = "  +0.03%  "

$ws.Range("E26").Value = "  +0.18%  "

$ws.Range("E27").Value = "  +3.13%  "

$ws.Range("D28").Value = "'9.14"
$ws.Range("E28").Value = "  +1.90%  "

$ws.Range("E29").Value = "  +2.85%  "

$ws.Range("D30").Value = "'2.85"
$ws.Range("E30").Value = "  +3.11%  "

$ws.Range("D31").Value = "'6.80"
$ws.Range("E31").Value = "  +7.38%  "

$ws.Range("D32").Value = "'28.07"
$ws.Range("E32").Value = "  +0.87%  "

$ws.Range("D33").Value = "'1.22"
$ws.Range("E33").Value = "  +1.58%  "

$ws.Range("E34").Value = "  -0.01%  "

$ws.Range("D35").Value = "'6.59"
$ws.Range("E35").Value = "  +0.48%  "

$ws.Range("D36").Value = "'55.41"
$ws.Range("E36").Value = "  -0.51%  "

$ws.Range("E37").Value = "  +1.02%  "

$ws.Range("D38").Value = "'480.42"
$ws.Range("E38").Value = "  +1.17%  "

$ws.Range("D39").Value = "'0.0420"
$ws.Range("E39").Value = "  +0.00%  "

$ws.Range("D40").Value = "'2.95"
$ws.Range("E40").Value = "  -6.52%  "

$ws.Range("D41").Value = "'8.87"
$ws.Range("E41").Value = "  +2.63%  "

$ws.Range("E42").Value = "  +3.92%  "

$ws.Range("D43").Value = "'0.119"
$ws.Range("E43").Value = "  +0.65%  "

$ws.Range("D44").Value = "2.938.43"
$ws.Range("E44").Value = "  -4.05%  "

$ws.Range("E45").Value = "  +1.46%  "

$ws.Range("D46").Value = "0.0₃0642"
$ws.Range("E46").Value = "  +5.01%  "

$ws.Range("D47").Value = "'28.61"
$ws.Range("E47").Value = "  -1.57%  "

$ws.Range("E48").Value = "  +0.07%  "

$ws.Range("E49").Value = "  +0.70%  "

$ws.Range("E50").Value = "  +2.42%  "

$ws.Range("D51").Value = "'120.58"
$ws.Range("E51").Value = "  +0.02%  "

